# Update the NATMI LR-pair (Slit2-Robo2) results with newly recomputed
# TPM-based values. Only the numeric data cells (columns E:T, rows 2:10)
# change; headers, labels and formatting are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.04558966666666667
$ws.Range("H2").Value = 0.136769
$ws.Range("I2").Value = 0.02375599288687187
$ws.Range("J2").Value = 0.02375599288687187
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.892869333333334
$ws.Range("N2").Value = 8.678608000000001
$ws.Range("O2").Value = 0.9698236995656885
$ws.Range("P2").Value = 0.9698236995656884
$ws.Range("Q2").Value = 0.1318849486168889
$ws.Range("R2").Value = 1.186964537552
$ws.Range("S2").Value = 0.02303912490840226
$ws.Range("T2").Value = 0.02303912490840225
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.04558966666666667
$ws.Range("H3").Value = 0.136769
$ws.Range("I3").Value = 0.02375599288687187
$ws.Range("J3").Value = 0.02375599288687187
$ws.Range("O3").Value = 0.006810193051573731
$ws.Range("P3").Value = 0.00681019305157373
$ws.Range("Q3").Value = 0.0009261084886666667
$ws.Range("R3").Value = 0.008334976398
$ws.Range("S3").Value = 0.0001617828976914098
$ws.Range("T3").Value = 0.0001617828976914098
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.04558966666666667
$ws.Range("H4").Value = 0.136769
$ws.Range("I4").Value = 0.02375599288687187
$ws.Range("J4").Value = 0.02375599288687187
$ws.Range("O4").Value = 0.02336610738273784
$ws.Range("P4").Value = 0.02336610738273783
$ws.Range("Q4").Value = 0.003177523783888889
$ws.Range("R4").Value = 0.028597714055
$ws.Range("S4").Value = 0.0005550850807782043
$ws.Range("T4").Value = 0.0005550850807782042
$ws.Range("I5").Value = 0.1978186777627204
$ws.Range("J5").Value = 0.1978186777627204
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.892869333333334
$ws.Range("N5").Value = 8.678608000000001
$ws.Range("O5").Value = 0.9698236995656885
$ws.Range("P5").Value = 0.9698236995656884
$ws.Range("Q5").Value = 1.098219985013333
$ws.Range("R5").Value = 9.883979865120001
$ws.Range("S5").Value = 0.1918492419110343
$ws.Range("T5").Value = 0.1918492419110343
$ws.Range("I6").Value = 0.1978186777627204
$ws.Range("J6").Value = 0.1978186777627204
$ws.Range("O6").Value = 0.006810193051573731
$ws.Range("P6").Value = 0.00681019305157373
$ws.Range("Q6").Value = 0.00771180382
$ws.Range("R6").Value = 0.06940623438
$ws.Range("S6").Value = 0.001347183384771181
$ws.Range("T6").Value = 0.001347183384771181
$ws.Range("I7").Value = 0.1978186777627204
$ws.Range("J7").Value = 0.1978186777627204
$ws.Range("O7").Value = 0.02336610738273784
$ws.Range("P7").Value = 0.02336610738273783
$ws.Range("S7").Value = 0.004622252466914938
$ws.Range("T7").Value = 0.004622252466914937
$ws.Range("I8").Value = 0.7784253293504076
$ws.Range("J8").Value = 0.7784253293504078
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.892869333333334
$ws.Range("N8").Value = 8.678608000000001
$ws.Range("O8").Value = 0.9698236995656885
$ws.Range("P8").Value = 0.9698236995656884
$ws.Range("Q8").Value = 4.321544675162667
$ws.Range("R8").Value = 38.89390207646401
$ws.Range("S8").Value = 0.7549353327462518
$ws.Range("T8").Value = 0.7549353327462519
$ws.Range("I9").Value = 0.7784253293504076
$ws.Range("J9").Value = 0.7784253293504078
$ws.Range("O9").Value = 0.006810193051573731
$ws.Range("P9").Value = 0.00681019305157373
$ws.Range("Q9").Value = 0.03034629235400001
$ws.Range("S9").Value = 0.005301226769111139
$ws.Range("T9").Value = 0.005301226769111139
$ws.Range("I10").Value = 0.7784253293504076
$ws.Range("J10").Value = 0.7784253293504078
$ws.Range("O10").Value = 0.02336610738273784
$ws.Range("P10").Value = 0.02336610738273783
$ws.Range("Q10").Value = 0.01818876983504469
$ws.Range("R10").Value = 0.9370765973850002
$ws.Range("S10").Value = 0.01818876983504469
$ws.Range("T10").Value = 0.01818876983504469
